# ============================================================================
# Scheduled runner update: refresh Universalis market-price derived columns
# (currentAveragePrice / currentAveragePriceNQ / currentAveragePriceHQ /
#  LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ) on the Leve Profit
# tables of every crafting-class worksheet.
# ============================================================================

$wb = $excel.ActiveWorkbook

# --- ALC!28 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1136.6666
$ws.Range("I28").Value = 1428.4
$ws.Range("J28").Value = 871.4545000000001
$ws.Range("K28").Value = 1428.4
$ws.Range("L28").Value = 871.4545000000001
$ws.Range("M28").Value = -943.4000000000001
$ws.Range("N28").Value = -1841.4545

# --- ALC!32 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 3173.3125
$ws.Range("J32").Value = 2534.818
$ws.Range("L32").Value = 2534.818
$ws.Range("N32").Value = -3186.818

# --- ALC!62 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 4576
$ws.Range("I62").Value = 3500
$ws.Range("J62").Value = 4647.7334
$ws.Range("K62").Value = 3500
$ws.Range("L62").Value = 4647.7334
$ws.Range("M62").Value = -2876
$ws.Range("N62").Value = -5895.7334

# --- ALC!65 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 4576
$ws.Range("I65").Value = 3500
$ws.Range("J65").Value = 4647.7334
$ws.Range("K65").Value = 17500
$ws.Range("L65").Value = 23238.667
$ws.Range("M65").Value = -14380
$ws.Range("N65").Value = -29478.667

# --- ALC!74 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 3750
$ws.Range("I74").Value = 3000
$ws.Range("J74").Value = 4500
$ws.Range("K74").Value = 3000
$ws.Range("L74").Value = 4500
$ws.Range("M74").Value = -2064
$ws.Range("N74").Value = -6372

# --- ALC!75 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H75").Value = 40000
$ws.Range("J75").Value = 40000
$ws.Range("L75").Value = 40000
$ws.Range("N75").Value = -41872

# --- ALC!77 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H77").Value = 3750
$ws.Range("I77").Value = 3000
$ws.Range("J77").Value = 4500
$ws.Range("K77").Value = 15000
$ws.Range("L77").Value = 22500
$ws.Range("M77").Value = -10320
$ws.Range("N77").Value = -31860

# --- ALC!78 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H78").Value = 40000
$ws.Range("J78").Value = 40000
$ws.Range("L78").Value = 120000
$ws.Range("N78").Value = -129360

# --- ALC!86 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 4312
$ws.Range("I86").Value = 2650
$ws.Range("J86").Value = 4836.8423
$ws.Range("K86").Value = 2650
$ws.Range("L86").Value = 4836.8423
$ws.Range("M86").Value = -1527
$ws.Range("N86").Value = -7082.8423

# --- ALC!88 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 647
$ws.Range("I88").Value = 695
$ws.Range("J88").Value = 615
$ws.Range("K88").Value = 695
$ws.Range("L88").Value = 615
$ws.Range("M88").Value = -289
$ws.Range("N88").Value = -1427

# --- ALC!89 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 4312
$ws.Range("I89").Value = 2650
$ws.Range("J89").Value = 4836.8423
$ws.Range("K89").Value = 13250
$ws.Range("L89").Value = 24184.2115
$ws.Range("M89").Value = -7634
$ws.Range("N89").Value = -35416.2115

# --- ALC!91 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H91").Value = 647
$ws.Range("I91").Value = 695
$ws.Range("J91").Value = 615
$ws.Range("K91").Value = 695
$ws.Range("L91").Value = 615
$ws.Range("M91").Value = 709
$ws.Range("N91").Value = -3423

# --- ALC!100 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 1560.8695
$ws.Range("I100").Value = 1666.6666
$ws.Range("J100").Value = 1445.4546
$ws.Range("K100").Value = 1666.6666
$ws.Range("L100").Value = 1445.4546
$ws.Range("M100").Value = -1125.6666
$ws.Range("N100").Value = -2527.4546

# --- ARM!2 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1132767.6
$ws.Range("J2").Value = 4903949
$ws.Range("L2").Value = 4903949
$ws.Range("N2").Value = -4904175

# --- ARM!9 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H9").Value = 15666.667
$ws.Range("J9").Value = 15666.667
$ws.Range("L9").Value = 15666.667
$ws.Range("N9").Value = -16006.667

# --- ARM!20 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H20").Value = 15666.667
$ws.Range("J20").Value = 15666.667
$ws.Range("L20").Value = 15666.667
$ws.Range("N20").Value = -16206.667

# --- ARM!116 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 1132767.6
$ws.Range("J116").Value = 4903949
$ws.Range("L116").Value = 4903949
$ws.Range("N116").Value = -4908537

# --- BSM!3 ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1132767.6
$ws.Range("J3").Value = 4903949
$ws.Range("L3").Value = 4903949
$ws.Range("N3").Value = -4904177

# --- CRP!16 ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1592.05
$ws.Range("I16").Value = 1213
$ws.Range("J16").Value = 1796.1538
$ws.Range("K16").Value = 1213
$ws.Range("L16").Value = 1796.1538
$ws.Range("M16").Value = -926
$ws.Range("N16").Value = -2370.1538

# --- CRP!60 ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H60").Value = 9363.125
$ws.Range("J60").Value = 10484.167
$ws.Range("L60").Value = 10484.167
$ws.Range("N60").Value = -11506.167

# --- CRP!68 ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H68").Value = 17550
$ws.Range("J68").Value = 19960
$ws.Range("L68").Value = 19960
$ws.Range("N68").Value = -21458

# --- CRP!71 ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H71").Value = 17550
$ws.Range("J71").Value = 19960
$ws.Range("L71").Value = 59880
$ws.Range("N71").Value = -67368

# --- CRP!113 ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 1592.05
$ws.Range("I113").Value = 1213
$ws.Range("J113").Value = 1796.1538
$ws.Range("K113").Value = 1213
$ws.Range("L113").Value = 1796.1538
$ws.Range("M113").Value = 957
$ws.Range("N113").Value = -6136.1538

# --- CUL!131 ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1976.4884
$ws.Range("I131").Value = 6806.25
$ws.Range("J131").Value = 1481.1282
$ws.Range("K131").Value = 20418.75
$ws.Range("L131").Value = 4443.3846
$ws.Range("M131").Value = -15378.75
$ws.Range("N131").Value = -14523.3846

# --- GSM!97 ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 2004.762
$ws.Range("I97").Value = 1982
$ws.Range("J97").Value = 2061.6667
$ws.Range("K97").Value = 1982
$ws.Range("L97").Value = 2061.6667
$ws.Range("M97").Value = -1486
$ws.Range("N97").Value = -3053.6667

# --- GSM!107 ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 563.4583
$ws.Range("I107").Value = 338.4
$ws.Range("J107").Value = 724.2143
$ws.Range("K107").Value = 338.4
$ws.Range("L107").Value = 724.2143
$ws.Range("M107").Value = 1581.6
$ws.Range("N107").Value = -4564.2143

# --- GSM!132 ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2021.9117
$ws.Range("I132").Value = 1830.0435
$ws.Range("J132").Value = 2423.0908
$ws.Range("K132").Value = 5490.1305
$ws.Range("L132").Value = 7269.2724
$ws.Range("M132").Value = -2960.1305
$ws.Range("N132").Value = -12329.2724

# --- LTW!23 ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H23").Value = 2983.3333
$ws.Range("I23").Value = 2983.3333
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 2983.3333
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = -2753.3333
$ws.Range("N23").ClearContents()

# --- LTW!93 ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1343.2
$ws.Range("I93").Value = 1343.2
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 1343.2
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = -95.20000000000005
$ws.Range("N93").ClearContents()

# --- WVR!113 ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 577.15
$ws.Range("I113").Value = 491.36365
$ws.Range("J113").Value = 682
$ws.Range("K113").Value = 1474.09095
$ws.Range("L113").Value = 2046
$ws.Range("M113").Value = 695.90905
$ws.Range("N113").Value = -6386
